# Apply Aug 13 2023 cryptos-list refresh: updated prices / 1h-volume%
# for most rows, plus a 3-row reorder (Algorand / RenderToken / EnergySwap)
# at rows 47-49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.386.78"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.849.45"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'240.64"
$ws.Range("D6").Value = "'0.6298"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07612"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.2916"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "'24.64"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'0.07758"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'5.024"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'0.6813"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "'83.15"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "'6.127"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "29.408.56"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "'229.47"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'12.34"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'7.477"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'159.32"
$ws.Range("D24").Value = "'0.1393"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'8.455"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'1.445"
$ws.Range("E27").Value = "  +9.41%  "
$ws.Range("D28").Value = "'1.472"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'0.05633"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "'4.114"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'4.052"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "'0.6968"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "'2.585"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'0.01830"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").Value = "1.238.55"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "'2.727"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").Value = "'6.425"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Value = "'0.9030"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'101.48"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'7.153"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'0.00000000117"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "'0.4004"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1156"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.686"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.992"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'0.05701"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "'0.4625"
$ws.Range("E51").Value = "  -0.21%  "
